# Updated test data for DC, TripCurrent, Voltdrop, BatteryStandby
#
# Sheet 1 = "Add_Devices_LoopA", Sheet 2 = "Add_Devices_LoopB"
#
# For both sheets:
#   - the User Story cell (B4) text changes from
#     "NGC-1826/TC-63775" to "NGC-1826/T916 OR TC-63775" and loses its
#     border/fill formatting (reset to the default/general style).
#   - the "worst case current" values in G1/G2 are rounded to whole numbers.
#   - two new columns are added (H = "Loop", I = "Column") describing the
#     built-in loop/column used for the DC unit calculation, with four
#     supporting rows (Built-in Loop-A..D) taking over the plain bordered
#     style that B4 used to have.
#   - the active selection moves to B4.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$sheetIndexes = @(1, 2)
$g1Values = @(300, 320)
$g2Values = @(220, 300)

for ($i = 0; $i -lt $sheetIndexes.Length; $i++) {

    $ws = $wb.Worksheets.Item($sheetIndexes[$i])

    # --- new "Loop" / "Column" header cells (H1:I1), formatted like E1 ---
    $ws.Range("E1").Copy()
    $ws.Range("H1").PasteSpecial($xlPasteFormats)
    $ws.Range("I1").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("H1").Value = "Loop"
    $ws.Range("I1").Value = "Column"

    # --- supporting Built-in Loop rows, formatted like the old B4 cell ---
    $ws.Range("B4").Copy()
    $ws.Range("H2").PasteSpecial($xlPasteFormats)
    $ws.Range("I2").PasteSpecial($xlPasteFormats)
    $ws.Range("H3").PasteSpecial($xlPasteFormats)
    $ws.Range("H4").PasteSpecial($xlPasteFormats)
    $ws.Range("H5").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("H2").Value = "Built-in Loop-A"
    $ws.Range("I2").Value = 2
    $ws.Range("H3").Value = "Built-in Loop-B"
    $ws.Range("H4").Value = "Built-in Loop-C"
    $ws.Range("H5").Value = "Built-in Loop-D"

    # --- update the User Story reference and drop its old formatting ---
    $ws.Range("B4").ClearFormats()
    $ws.Range("B4").Value = "NGC-1826/T916 OR TC-63775"

    # --- rounded worst-case current values ---
    $ws.Range("G1").Value = $g1Values[$i]
    $ws.Range("G2").Value = $g2Values[$i]

    # --- widen column H slightly to fit the new "Built-in Loop-X" text ---
    $ws.Columns.Item(8).ColumnWidth = 12.5

    # --- move the selection to B4, matching the latest edit location ---
    $ws.Range("B4").Select()
}

# Restore the originally active tab (Add_Devices_LoopA) since selecting
# ranges on the second sheet above shifts the active tab to it.
$wb.Worksheets.Item(1).Activate()

Write-Output "Applied DC unit loading updates to both loop sheets."
